# Insert a new data row at row 280, shifting all subsequent rows down by one
# (the final row's content is preserved by the shift, extending the used
# range from R372 to R373), then populate the new row 280 with the new
# record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("280:280").Insert()

$ws.Cells.Item(280, 1).Value = 5
$ws.Cells.Item(280, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(280, 3).Value = "Maule"
$ws.Cells.Item(280, 4).Value = 44559
$ws.Cells.Item(280, 5).Value = 7
$ws.Cells.Item(280, 6).Value = 100114001
$ws.Cells.Item(280, 7).Value = "Papa"
$ws.Cells.Item(280, 8).Value = "Asterix"
$ws.Cells.Item(280, 9).Value = "1a nueva(o)"
$ws.Cells.Item(280, 10).Value = 2000
$ws.Cells.Item(280, 11).Value = 7000
$ws.Cells.Item(280, 12).Value = 8000
$ws.Cells.Item(280, 13).Value = 7600
$ws.Cells.Item(280, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(280, 15).Value = "Región del Maule"
$ws.Cells.Item(280, 16).Value = 304
$ws.Cells.Item(280, 17).Value = 25
$ws.Cells.Item(280, 18).Value = "Hortaliza"
